# Add three new vocabulary entries (fall, scatter, scratch) as new rows
# 100-102 at the bottom of the word list on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 100: fall
$ws.Cells.Item(100, 1).Value = "fall"
$ws.Cells.Item(100, 2).Value = "to move or drop down from a higher position to a lower position"
$ws.Cells.Item(100, 3).Value = "The book fell from his hands."
$ws.Cells.Item(100, 4).Value = "the tree falls."
$ws.Rows.Item(100).RowHeight = 45

# Row 101: scatter
$ws.Cells.Item(101, 1).Value = "scatter"
$ws.Cells.Item(101, 2).Value = "if someone scatters a lot of things, or if they scatter, they are thrown or dropped over a wide area in an irregular way"
$ws.Cells.Item(101, 3).Value = "the wind scattered the dry fallen leaves."
$ws.Cells.Item(101, 4).Value = "Scatter the onions over the fish."
$ws.Rows.Item(101).RowHeight = 75

# Row 102: scratch
$ws.Cells.Item(102, 1).Value = "scratch"
$ws.Cells.Item(102, 2).Value = "to rub your skin with your nails because it feels uncomfortable"
$ws.Cells.Item(102, 3).Value = "the cat scratched the living room carpet."
$ws.Cells.Item(102, 4).Value = "John yawned and scratched his leg."
$ws.Rows.Item(102).RowHeight = 45

# Update the active selection to match the new editing position.
$null = $ws.Range("F100").Select()
